$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Recorded By" (column G) email-list re-orderings (same underlying set of
# recipients, just re-sequenced) plus a couple of genuine roster additions.
# ---------------------------------------------------------------------------

$group_servinaz = "alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G2").Value = $group_servinaz
$ws.Range("G21").Value = $group_servinaz
$ws.Range("G40").Value = $group_servinaz

$group_sara = "Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G6").Value = $group_sara
$ws.Range("G44").Value = $group_sara

$group_wafaa = "wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G18").Value = $group_wafaa
$ws.Range("G37").Value = $group_wafaa
$ws.Range("G56").Value = $group_wafaa
$ws.Range("G75").Value = $group_wafaa
$ws.Range("G94").Value = $group_wafaa
$ws.Range("G113").Value = $group_wafaa

$group_neveen5 = "marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G19").Value = $group_neveen5
$ws.Range("G76").Value = $group_neveen5
$ws.Range("G95").Value = $group_neveen5

$group_neveen6 = "nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G20").Value = $group_neveen6
$ws.Range("G38").Value = $group_neveen6
$ws.Range("G39").Value = $group_neveen6
$ws.Range("G57").Value = $group_neveen6
$ws.Range("G58").Value = $group_neveen6
$ws.Range("G77").Value = $group_neveen6
$ws.Range("G96").Value = $group_neveen6
$ws.Range("G115").Value = $group_neveen6

$group_lamiaa3 = "AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg"
$ws.Range("G24").Value = $group_lamiaa3
$ws.Range("G81").Value = $group_lamiaa3

$group_norhan = "yasmintarek@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G25").Value = $group_norhan
$ws.Range("G82").Value = $group_norhan

$ws.Range("G28").Value = "dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"

$group_lamiaa6 = "lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("G43").Value = $group_lamiaa6
$ws.Range("G100").Value = $group_lamiaa6

$group_mohammed = "Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G59").Value = $group_mohammed
$ws.Range("G78").Value = $group_mohammed
$ws.Range("G97").Value = $group_mohammed

$group_dina = "dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G66").Value = $group_dina
$ws.Range("G85").Value = $group_dina

# G47 gains two extra recipients (merna.said, arwaelsayed03, maryam.ahmed)
$ws.Range("G47").Value = "merna.said@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
$ws.Range("H47").Value = "92/216"

# ---------------------------------------------------------------------------
# Plain numeric statistic updates (Class Statistics boxes)
# ---------------------------------------------------------------------------
$ws.Range("L7").Value = 3
$ws.Range("L8").Value = 73
$ws.Range("P20").Value = 2
$ws.Range("Q20").Value = 12

# ---------------------------------------------------------------------------
# Percentage values stored as literal text (e.g. "44.9%"). Typing a
# percent-looking string directly gets auto-coerced to a numeric percentage
# by Excel's input parser, so force it in as text with a leading apostrophe
# and then re-stamp the original cell format (copy/paste-special formats
# from an unaffected sibling cell that still carries the source style) so
# the cell's style index is not disturbed by the quote-prefix flag.
# ---------------------------------------------------------------------------
$ws.Range("L10").Value = "'44.9%"
$ws.Range("L9").Copy()
$ws.Range("L10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("S17").Value = "'58.0%"
$ws.Range("S18").Copy()
$ws.Range("S17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row 104 (MICROBIOLOGY B3 session 1) flips from "Pending" (yellow) to
# "Not Recorded" (pink) status - copy the pink format from an existing
# "Not Recorded" row (62) onto row 104, then fix up the status text.
# ---------------------------------------------------------------------------
$ws.Range("A62:I62").Copy()
$ws.Range("A104:I104").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("I104").Value = "Not Recorded"
